# Monthly rollover update.
#
# Sheet "VENTAS POR GRUPO": a batch of per-client/product cells that still
# carried a stray prior-period amount get cleared back to 0, and the row 60
# "<n> de 58" running-count labels that summarised those now-cleared columns
# are reset to "0 de 58".
#
# Sheet "VENTA MENSUAL": the rolling 4-month window (agosto..noviembre)
# advances by one month (septiembre..diciembre) - every row's data shifts
# one column to the left (C<-D, D<-E, E<-F) and the newly opened rightmost
# month column F is zeroed out, ready to accumulate. Column widths for the
# three shifted columns follow the same shift. PRESUPUESTO (G) is untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$zeroRefs = @(
    "M2", "K4", "M7",
    "E10", "I10", "L10", "M10",
    "L16", "M16",
    "M27",
    "D29", "L29",
    "D30", "L30", "M30",
    "D35", "M35",
    "M46",
    "M47", "P47",
    "M52",
    "M59"
)
foreach ($ref in $zeroRefs) {
    $ws1.Range($ref).Value2 = 0
}

$resetCountRefs = @("D60", "E60", "I60", "K60", "L60", "M60", "P60")
foreach ($ref in $resetCountRefs) {
    $ws1.Range($ref).Value2 = "0 de 58"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Shift each data row's agosto/sept/oct/nov values left by one month;
# the new rightmost month (diciembre) starts at 0. Read the old D/E/F
# before overwriting so the shift is correct row by row.
for ($r = 2; $r -le 60; $r++) {
    $oldD = $ws2.Cells.Item($r, 4).Value2
    $oldE = $ws2.Cells.Item($r, 5).Value2
    $oldF = $ws2.Cells.Item($r, 6).Value2

    $ws2.Cells.Item($r, 3).Value2 = $oldD
    $ws2.Cells.Item($r, 4).Value2 = $oldE
    $ws2.Cells.Item($r, 5).Value2 = $oldF
    $ws2.Cells.Item($r, 6).Value2 = 0
}

# Header labels follow the same one-month advance.
$ws2.Range("C1").Value2 = "septiembre"
$ws2.Range("D1").Value2 = "octubre"
$ws2.Range("E1").Value2 = "noviembre"
$ws2.Range("F1").Value2 = "diciembre"

# Column widths follow the same shift (C<-old D width, D<-old E width,
# E<-old F width). ColumnWidth is offset from the stored sheet width by a
# fixed padding of 5/6 character, so subtract that to land exactly on the
# target stored widths of 16 / 14 / 15.
$ws2.Columns.Item(3).ColumnWidth = 16 - 0.8333333333333334
$ws2.Columns.Item(4).ColumnWidth = 14 - 0.8333333333333334
$ws2.Columns.Item(5).ColumnWidth = 15 - 0.8333333333333334
